# "Capitulo 2 - Alteração"
# Adds a new slide 2 ("Somente Título" / Title Only layout) right after the
# existing title slide, with a bold "CAP 2 " title that hyperlinks out to the
# external "Cap 2.docx" file - mirroring the "Cap 1" link already used on
# slide 1.

$p = $ppt.ActivePresentation

# ppLayoutTitleOnly = 11 -> maps to slideLayout6.xml ("Somente Título")
$s = $p.Slides.Add(2, 11)

$sh = $s.Shapes.Item(1)
$sh.Name = "Título 1"

$tf = $sh.TextFrame
$tf.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$tr = $tf.TextRange
$tr.Text = "CAP"
$tr.Font.Size = 28
$tr.Font.Bold = $true
$tr.ActionSettings.Item(1).Hyperlink.Address = "file:///C:\Users\Aluno\Desktop\Jovana%20TI%20-%20Noite\Cap%202.docx"

$tr2 = $tr.InsertAfter(" 2 ")
$tr2.Font.Size = 28
$tr2.Font.Bold = $true
$tr2.ActionSettings.Item(1).Hyperlink.Address = "file:///C:\Users\Aluno\Desktop\Jovana%20TI%20-%20Noite\Cap%202.docx"
